$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.783.39"
$ws.Range("E2").Value = "  -2.76%  "

$ws.Range("D3").Value = "3.195.17"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "599.50"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "153.15"
$ws.Range("E6").Value = "  -3.51%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "3.193.13"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -3.59%  "

$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  -4.36%  "

$ws.Range("E11").Value = "  -1.75%  "

$ws.Range("D12").Value = "0.480"
$ws.Range("E12").Value = "  -5.71%  "

$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -5.01%  "

$ws.Range("D14").Value = "37.28"
$ws.Range("E14").Value = "  -4.56%  "

$ws.Range("D15").Value = "3.722.73"
$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").Value = "64.883.22"
$ws.Range("E16").Value = "  -2.60%  "

$ws.Range("D17").Value = "3.202.55"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "7.09"
$ws.Range("E19").Value = "  -5.32%  "

$ws.Range("D20").Value = "485.85"
$ws.Range("E20").Value = "  -5.19%  "

$ws.Range("D21").Value = "14.88"
$ws.Range("E21").Value = "  -2.92%  "

$ws.Range("D22").Value = "0.722"
$ws.Range("E22").Value = "  -1.94%  "

$ws.Range("D23").Value = "7.82"
$ws.Range("E23").Value = "  -3.23%  "

$ws.Range("D24").Value = "13.98"
$ws.Range("E24").Value = "  -5.73%  "

$ws.Range("D25").Value = "85.85"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").Value = "2.94"
$ws.Range("E27").Value = "  -1.77%  "

$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -5.45%  "

$ws.Range("D29").Value = "0.128"
$ws.Range("E29").Value = "  +34.52%  "

$ws.Range("D30").Value = "2.29"
$ws.Range("E30").Value = "  -5.03%  "

$ws.Range("D31").Value = "6.98"
$ws.Range("E31").Value = "  -1.03%  "

$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  -9.18%  "

$ws.Range("D33").Value = "27.09"
$ws.Range("E33").Value = "  -4.21%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("E35").Value = "  -6.37%  "

$ws.Range("D36").Value = "6.17"
$ws.Range("E36").Value = "  -5.86%  "

$ws.Range("E37").Value = "  +8.55%  "

$ws.Range("D38").Value = "54.66"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("D39").Value = "480.06"
$ws.Range("E39").Value = "  -7.04%  "

$ws.Range("D40").Value = "0.0₃0729"
$ws.Range("E40").Value = "  -6.10%  "

$ws.Range("D41").Value = "0.0407"
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("D42").Value = "0.125"
$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("D43").Value = "8.57"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  -1.66%  "

$ws.Range("D45").Value = "2.925.62"
$ws.Range("E45").Value = "  +1.87%  "

$ws.Range("D46").Value = "0.279"
$ws.Range("E46").Value = "  -7.23%  "

$ws.Range("D47").Value = "27.64"
$ws.Range("E47").Value = "  -3.45%  "

$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").Value = "120.77"
$ws.Range("E51").Value = "  -2.05%  "
